$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells AD1:AF1 with the same style as existing headers (row 1)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$headerRange = $ws.Range("AD1:AF1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108  # xlCenter
$headerRange.VerticalAlignment = -4160    # xlTop
$headerRange.Borders.LineStyle = 1
$headerRange.Borders.Weight = 2

# Fill in the Wins/Losses/Ties values for each data row (2-51)
for ($r = 2; $r -le 51; $r++) {
    $ws.Cells.Item($r, 30).Value = 85  # AD
    $ws.Cells.Item($r, 31).Value = 76  # AE
    $ws.Cells.Item($r, 32).Value = 0   # AF
}
